$d = $word.ActiveDocument

# Paragraph 4 is the bullet that starts with "Contexts / ContextKinds: ...".
# It gets split into 8 bullets total:
#   1-6: brand new bullets (clean paragraph-mark formatting, like most other
#        bullets in the doc)
#   7:   the original bullet, text replaced with the "Rule: (Amor, ...)" line
#        (keeps its original paragraph-mark formatting, which carries an
#        explicit "no underline" override)
#   8:   a brand new bullet ("Statements: (unAmor, ...)") that inherits that
#        same "no underline" paragraph-mark formatting from #7.

$newTexts = @(
  "Contexts / ContextKinds: Implications / Rules (Upper asserted / Aligned Knowledge / Primitives). (Context, LHS, Concept, RHS);",
  "Implication / Assertion Statement: (Context, Subject, Property, ModelObject);",
  "ContextKind Aggregation. Instance: Context, Atribute: Concept, Value: RHS, from Assertion Statements.",
  "Rule Statement: (Context: ContextKind, LHS: SubjectKind, Concept: PropertyKind. RHS: ObjectKind.",
  "Aggregate Kinds into Rule Statements. (KindStatements). Statements match Contexts, match SK, apply Concept, match OK.",
  "Materialize Rule Statements from Model into CSPOs Occurrences from Resources:"
)

$target = $d.Paragraphs.Item(4)

# Insert 6 new (clean-formatted) paragraphs right after paragraph 3, i.e.
# right before the target paragraph, one at a time, filling each with its
# text as we go so paragraph indices stay predictable.
$anchor = $d.Paragraphs.Item(3)
for ($i = 0; $i -lt $newTexts.Length; $i++) {
  $anchor.Range.InsertParagraphAfter()
  $created = $d.Paragraphs.Item(4 + $i)
  $created.Range.Text = $newTexts[$i]
  $anchor = $created
}

# The original paragraph (now shifted down to slot 10) gets its text replaced
# in place, preserving its original paragraph-mark / run formatting.
$ruleP = $d.Paragraphs.Item(4 + $newTexts.Length)
$ruleP.Range.Text = "Rule: (Amor, Amante, Ama, Amada);"

# Append the final new bullet right after it, inheriting its formatting.
$ruleP = $d.Paragraphs.Item(4 + $newTexts.Length)
$ruleP.Range.InsertParagraphAfter()
$lastP = $d.Paragraphs.Item(4 + $newTexts.Length + 1)
$lastP.Range.Text = "Statements: (unAmor, pedro, amaA, maría);"

Write-Output "done"
